$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The report table (rows 7-10) lists short-stock items alphabetically.
# Two new items were added to the source data:
#   - "ASPIRIN PROTECT 100MG 30 GASTRO-RESISTANT TAB" (goes before EPIMAG)
#   - "PANADOL ADVANCE 500 MG 48 TABLETS" (goes before SPINOBAC)
# Insert two blank rows at the right spots, shifting everything else (incl.
# the totals row and the footer row) down, then fill the new rows in with
# the same look &amp; feel (styles / merges / row height) as their neighbours.
# ---------------------------------------------------------------------------

# Row 7: blank row that will hold the new ASPIRIN line (pushes EPIMAG -> 8,
# GYNERA -> 9, SPINOBAC -> 10, Calona -> 11, totals -> 12, footer -> 13).
$ws.Rows("7:7").Insert()

# Row 10: after the row above shifted SPINOBAC down to row 10, insert another
# blank row there for PANADOL (pushes SPINOBAC -> 11, Calona -> 12,
# totals -> 13, footer -> 14).
$ws.Rows("10:10").Insert()

# --- Copy formatting (styles, fonts, number formats, fills, borders) from
# the row right below each freshly inserted blank row so the new rows match
# the rest of the table exactly. ---
$ws.Range("A8:Q8").Copy()
$ws.Range("A7:Q7").PasteSpecial(-4122)

$ws.Range("A11:Q11").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Recreate the merged cells for the two new data rows (same pattern as
# every other data row: A:B, C:G, H:K, L:M, N:O). ---
$ws.Range("A7:B7").Merge()
$ws.Range("C7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("N7:O7").Merge()

$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

# --- Row heights matching the source rows they were copied from. ---
$ws.Rows("7:7").RowHeight = 25.5
$ws.Rows("10:10").RowHeight = 25.5

# --- Fill in the new ASPIRIN row (row 7). Columns C/H/L/N/P/Q are stored as
# text in this report even when the content looks numeric, so force text
# before writing the values; L and P normally carry a numeric-looking number
# format (copied above), so stash it, write as plain text "@", then restore
# it afterwards so the cell keeps its original format but stays text-typed.
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "ASPIRIN PROTECT 100MG 30 GASTRO-RESISTANT TAB"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "4:0"
$fmtL7 = $ws.Range("L7").NumberFormat
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "1"
$ws.Range("L7").NumberFormat = $fmtL7
$ws.Range("N7").NumberFormat = "@"
$ws.Range("N7").Value = "78.00"
$fmtP7 = $ws.Range("P7").NumberFormat
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "78.0000"
$ws.Range("P7").NumberFormat = $fmtP7
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "1:0"

# --- Fill in the new PANADOL row (row 10). ---
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "PANADOL ADVANCE 500 MG 48 TABLETS"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "2:3"
$fmtL10 = $ws.Range("L10").NumberFormat
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "1"
$ws.Range("L10").NumberFormat = $fmtL10
$ws.Range("N10").NumberFormat = "@"
$ws.Range("N10").Value = "92.00"
$fmtP10 = $ws.Range("P10").NumberFormat
$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "23.0000"
$ws.Range("P10").NumberFormat = $fmtP10
$ws.Range("Q10").NumberFormat = "@"
$ws.Range("Q10").Value = "0:1"

# --- Renumber the "م" (row index) column for all six data rows now that two
# new rows were inserted. ---
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6

# --- Update the grand total (was 249, now 350) and the "printed at"
# timestamp in the footer (was 10:39 AM, now 10:45 AM). The totals row and
# the footer row shifted down from 11/12 to 13/14 after the two inserts. ---
$ws.Range("P13").Value = 350
$ws.Range("K14").Value = "Sunday, 14 September, 2025 10:45 AM"
